$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "filename",
    "approx_split",
    "Diad1_pos",
    "Diad2_pos",
    "HB1_pos",
    "HB2_pos",
    "C13_pos",
    "Diad1_abs_prom",
    "Diad2_abs_prom",
    "HB1_abs_prom",
    "HB2_abs_prom",
    "C13_abs_prom",
    "Mean_abs_HB_prom",
    "Diad2_HB2_abs_prom_ratio",
    "Diad1_HB1_abs_prom_ratio",
    "Diad1_rel_prom",
    "Diad2_rel_prom",
    "HB1_rel_prom",
    "HB2_rel_prom",
    "C13_rel_prom",
    "Diad1_HB1_abs_prom_ratio",
    "Diad2_HB2_abs_prom_ratio",
    "Diad1_HB1_Valley_prom",
    "Diad2_HB2_abs_prom_ratio",
    "Mean_Diad_HB_Valley_prom",
    "Mean_abs_HB_prom",
    "Diad1_prom/std_betweendiads",
    "Diad2_prom/std_betweendiads",
    "Diad2_height",
    "HB2_height",
    "C13_height",
    "Diad1_height",
    "HB1_height",
    "Diad1_Median_Bck",
    "Diad2_Median_Bck",
    "C13_HB2_abs_prom_ratio",
    "Diad2_HB2_Valley_prom"
)

# Column B is the 2nd column; the header row extends from B1 through AL1 (37 entries)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Make sure the newly-added cells (beyond the original AD1) pick up the same
# bold/bordered header style ("s=1") that the rest of the header row uses.
$src = $ws.Range("B1")
$src.Copy()
$dst = $ws.Range("AE1:AL1")
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
